# Automatic update of files.
# Applies the row 5-8 corrections described by the commit diff:
#  - Row 5 (Dicranum flagellare obs) and Row 8 (Sarcodon squamosus obs) swap
#    identity/content, each picking up corrected Id/Taxonsorteringsordning and
#    rounded-integer Ost/Nord (Q/R) coordinates, and both drop their
#    Starttid/Sluttid (Z/AB) time-stamp cells.
#  - Row 6 and Row 7 keep their species content but swap "Id" (A) values,
#    receive corrected rounded-integer Ost/Nord coordinates, drop their
#    Starttid/Sluttid (Z/AB) cells, and the "Publik kommentar" (AC) note
#    moves from row 6 to row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 5: becomes the Motaggsvamp / Sarcodon squamosus observation
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 111934086
$ws.Range("B5").Value = 90689
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5966
$ws.Range("F5").Value = "Motaggsvamp"
$ws.Range("G5").Value = "Sarcodon squamosus"
$ws.Range("H5").Value = "(Schaeff.) Quél."
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("P5").Value = "Tallskogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q5").Value = 413681
$ws.Range("R5").Value = 6586805
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

# ---------------------------------------------------------------------
# Row 6: stays Flagellkvastmossa / Dicranum flagellare, new Id + coords
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 111934050
$ws.Range("Q6").Value = 413638
$ws.Range("R6").Value = 6587077
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").ClearContents()

# ---------------------------------------------------------------------
# Row 7: stays Flagellkvastmossa / Dicranum flagellare, new Id + coords,
# gains the "Rätt riklig längs stigen" public comment moved from row 6
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 111934059
$ws.Range("Q7").Value = 413640
$ws.Range("R7").Value = 6586794
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AC7").Value = "Rätt riklig längs stigen"

# ---------------------------------------------------------------------
# Row 8: becomes the Flagellkvastmossa / Dicranum flagellare observation
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 111934066
$ws.Range("B8").Value = 93289
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 2170
$ws.Range("F8").Value = "Flagellkvastmossa"
$ws.Range("G8").Value = "Dicranum flagellare"
$ws.Range("H8").Value = "Hedw."
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("K8").Value = "med groddkorn"
$ws.Range("L8").Value = ""
$ws.Range("P8").Value = "Skogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q8").Value = 413590
$ws.Range("R8").Value = 6586912
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()
